# Weekly update: insert the newest price record at the top of the data
# (row 371), shifting all existing records down by one row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 371; existing rows 371-395 shift to 372-396.
$ws.Rows.Item(371).Insert()

# Populate the newly inserted row with the latest week's data.
$ws.Cells.Item(371, 1).Value2 = 4
$ws.Cells.Item(371, 2).Value2 = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(371, 3).Value2 = "Los Lagos"
$ws.Cells.Item(371, 4).Value2 = 45013
$ws.Cells.Item(371, 5).Value2 = 10
$ws.Cells.Item(371, 6).Value2 = 100112021
$ws.Cells.Item(371, 7).Value2 = "Ají"
$ws.Cells.Item(371, 8).Value2 = "Inferno"
$ws.Cells.Item(371, 9).Value2 = "Primera"
$ws.Cells.Item(371, 10).Value2 = 180
$ws.Cells.Item(371, 11).Value2 = 22000
$ws.Cells.Item(371, 12).Value2 = 22000
$ws.Cells.Item(371, 13).Value2 = 22000
$ws.Cells.Item(371, 14).Value2 = "$/caja 10 kilos"
$ws.Cells.Item(371, 15).Value2 = "Región de Arica y Parinacota"
$ws.Cells.Item(371, 16).Value2 = 2200
$ws.Cells.Item(371, 17).Value2 = 10
$ws.Cells.Item(371, 18).Value2 = "Hortaliza"

Write-Output "done"
